$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary section updates ---
$ws.Range("E11").Value = 284024
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 7

# The detail table shrinks from 8 worker/period rows (16-23) to 7 (16-22).
# Row 22 becomes the new last row of the table and must carry the "closing"
# bottom-border formatting that the old last row (23) had, so copy that
# formatting up before touching any values.
$ws.Range("B23:J23").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Replace the worker/period data table (rows 16-22) with the updated data set ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143379586"
$ws.Range("D16").Value = "RICARDO JAVIER PEREZ BERNAL"
$ws.Range("E16").Value = "1707"
$ws.Range("F16").Value = 3935
$ws.Range("G16").Value = 737717

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73203743"
$ws.Range("D17").Value = "WILSON MANUEL RECUERO SEJIN"
$ws.Range("E17").Value = "1909"
$ws.Range("F17").Value = 48533
$ws.Range("G17").Value = 1300000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "19562446"
$ws.Range("D18").Value = "JOSE GREGORIO LOBATO ESMERAL"
$ws.Range("E18").Value = "2504"
$ws.Range("F18").Value = 3796
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "19562446"
$ws.Range("D19").Value = "JOSE GREGORIO LOBATO ESMERAL"
$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "19562446"
$ws.Range("D20").Value = "JOSE GREGORIO LOBATO ESMERAL"
$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "19562446"
$ws.Range("D21").Value = "JOSE GREGORIO LOBATO ESMERAL"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "19562446"
$ws.Range("D22").Value = "JOSE GREGORIO LOBATO ESMERAL"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

# Old row 23 (previously RICARDO JAVIER PEREZ BERNAL, now consolidated into row 16)
# is no longer needed. Deleting it shifts the footer/signature rows (28,29) up to
# (27,28) and keeps the dimension + merged cells in sync automatically.
$ws.Rows("23").Delete()
